$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells K1, L1 and copy header style/format from J1
$ws.Range("K1").Value = "Coef_Test"
$ws.Range("L1").Value = "VIP_Test"
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

# Update I/J (recomputed 95% CI strings) and add K/L (Coef_Test / VIP_Test) for rows 2-53
$ws.Cells.Item(2, 9).Value = '[-0.1007629  -0.04512894 -0.07477275]'
$ws.Cells.Item(2, 10).Value = '[1.35980648 2.72339097 2.09163488]'
$ws.Cells.Item(2, 11).Value = -0.06998639338294839
$ws.Cells.Item(2, 12).Value = 1.954274798075686
$ws.Cells.Item(3, 9).Value = '[-0.04694931  0.01917425 -0.00468743]'
$ws.Cells.Item(3, 10).Value = '[0.11699492 0.62005291 0.30754726]'
$ws.Cells.Item(3, 11).Value = 0.02563189554066063
$ws.Cells.Item(3, 12).Value = 0.8984445610919483
$ws.Cells.Item(4, 9).Value = '[0.02179637 0.06491312 0.04797865]'
$ws.Cells.Item(4, 10).Value = '[0.6302717  1.94833668 1.36355491]'
$ws.Cells.Item(4, 11).Value = 0.02327437662045344
$ws.Cells.Item(4, 12).Value = 0.7602436885117888
$ws.Cells.Item(5, 9).Value = '[-0.07107525 -0.04031215 -0.05798911]'
$ws.Cells.Item(5, 10).Value = '[1.12898745 1.98087169 1.62055228]'
$ws.Cells.Item(5, 11).Value = -0.06035640889745132
$ws.Cells.Item(5, 12).Value = 1.716371401082523
$ws.Cells.Item(6, 9).Value = '[-0.0253416   0.04117707  0.00261577]'
$ws.Cells.Item(6, 10).Value = '[0.05577179 0.22376844 0.05633274]'
$ws.Cells.Item(6, 11).Value = 0.01671235831834629
$ws.Cells.Item(6, 12).Value = 0.4859594348977094
$ws.Cells.Item(7, 9).Value = '[-0.07931125 -0.0125513  -0.04252455]'
$ws.Cells.Item(7, 10).Value = '[0.41912617 2.12186608 1.26038971]'
$ws.Cells.Item(7, 11).Value = -0.04263943410674419
$ws.Cells.Item(7, 12).Value = 1.207535720748384
$ws.Cells.Item(8, 9).Value = '[-0.03652473 -0.001601   -0.01933979]'
$ws.Cells.Item(8, 10).Value = '[0.1815231  1.09443784 0.63310178]'
$ws.Cells.Item(8, 11).Value = -0.02055761007525408
$ws.Cells.Item(8, 12).Value = 0.6801933451040987
$ws.Cells.Item(9, 9).Value = '[-0.01737787  0.03395016  0.01492792]'
$ws.Cells.Item(9, 10).Value = '[0.27397484 1.17836577 0.61117948]'
$ws.Cells.Item(9, 11).Value = 0.01449337243158384
$ws.Cells.Item(9, 12).Value = 0.5642896212759561
$ws.Cells.Item(10, 9).Value = '[-0.00630157  0.0280316   0.0089574 ]'
$ws.Cells.Item(10, 10).Value = '[0.20854145 0.92656847 0.50538687]'
$ws.Cells.Item(10, 11).Value = 0.01262689918707682
$ws.Cells.Item(10, 12).Value = 0.6160078904820261
$ws.Cells.Item(11, 9).Value = '[-0.03374221  0.02510135 -0.00566471]'
$ws.Cells.Item(11, 10).Value = '[0.07864001 0.52502081 0.23869808]'
$ws.Cells.Item(11, 11).Value = 0.005132444511115117
$ws.Cells.Item(11, 12).Value = 0.3821630548723957
$ws.Cells.Item(12, 9).Value = '[-0.01718424  0.05133559  0.02946482]'
$ws.Cells.Item(12, 10).Value = '[0.22586699 1.35736275 0.85582104]'
$ws.Cells.Item(12, 11).Value = 0.02079697018988702
$ws.Cells.Item(12, 12).Value = 0.677901137445071
$ws.Cells.Item(13, 9).Value = '[-0.01420785  0.01878986  0.0011803 ]'
$ws.Cells.Item(13, 10).Value = '[0.15679531 0.59226728 0.36427733]'
$ws.Cells.Item(13, 11).Value = 0.01095419432883107
$ws.Cells.Item(13, 12).Value = 0.5789520736713877
$ws.Cells.Item(14, 9).Value = '[-0.0267616   0.01759392 -0.005316  ]'
$ws.Cells.Item(14, 10).Value = '[0.05838522 0.34358546 0.06260275]'
$ws.Cells.Item(14, 11).Value = -0.001543440509054234
$ws.Cells.Item(14, 12).Value = 0.3316123335420682
$ws.Cells.Item(15, 9).Value = '[-0.0300992  0.0116677 -0.0102094]'
$ws.Cells.Item(15, 10).Value = '[0.08776678 0.79318285 0.40930335]'
$ws.Cells.Item(15, 11).Value = -0.01527241278230754
$ws.Cells.Item(15, 12).Value = 0.582964702296681
$ws.Cells.Item(16, 9).Value = '[-0.09329397 -0.04530015 -0.07133276]'
$ws.Cells.Item(16, 10).Value = '[1.39750305 2.72755758 2.10592718]'
$ws.Cells.Item(16, 11).Value = -0.06762300956980039
$ws.Cells.Item(16, 12).Value = 1.888522438813666
$ws.Cells.Item(17, 9).Value = '[-0.05064637 -0.03102167 -0.04376295]'
$ws.Cells.Item(17, 10).Value = '[0.91731671 1.42420992 1.2060261 ]'
$ws.Cells.Item(17, 11).Value = -0.03473016232920018
$ws.Cells.Item(17, 12).Value = 1.056571374988954
$ws.Cells.Item(18, 9).Value = '[-0.0332044   0.02703819 -0.00069493]'
$ws.Cells.Item(18, 10).Value = '[0.13847268 0.28458354 0.13957726]'
$ws.Cells.Item(18, 11).Value = 0.02521006414153381
$ws.Cells.Item(18, 12).Value = 0.7942035389486203
$ws.Cells.Item(19, 9).Value = '[-0.04622381 -0.00428554 -0.02483723]'
$ws.Cells.Item(19, 10).Value = '[0.27887496 1.29413545 0.67387404]'
$ws.Cells.Item(19, 11).Value = -0.02302348207357939
$ws.Cells.Item(19, 12).Value = 0.6786725491352518
$ws.Cells.Item(20, 9).Value = '[-0.01613416  0.01014484 -0.00162836]'
$ws.Cells.Item(20, 10).Value = '[0.1627331  0.62254494 0.43245115]'
$ws.Cells.Item(20, 11).Value = -0.002349266397172134
$ws.Cells.Item(20, 12).Value = 0.5267533069281976
$ws.Cells.Item(21, 9).Value = '[-0.035987    0.02510434 -0.00071172]'
$ws.Cells.Item(21, 10).Value = '[0.08139516 0.52864701 0.22575644]'
$ws.Cells.Item(21, 11).Value = 0.0202346753620685
$ws.Cells.Item(21, 12).Value = 0.6586538111515085
$ws.Cells.Item(22, 9).Value = '[-0.02774111  0.02967525 -0.00566608]'
$ws.Cells.Item(22, 10).Value = '[0.0235926  0.62156374 0.30162123]'
$ws.Cells.Item(22, 11).Value = 0.007494858620690182
$ws.Cells.Item(22, 12).Value = 0.3977243013968144
$ws.Cells.Item(23, 9).Value = '[-0.04519474 -0.01390111 -0.02824776]'
$ws.Cells.Item(23, 10).Value = '[0.54614346 1.23286419 0.84080388]'
$ws.Cells.Item(23, 11).Value = -0.0170548403506805
$ws.Cells.Item(23, 12).Value = 0.664950450675809
$ws.Cells.Item(24, 9).Value = '[-0.00050222  0.03878739  0.02063863]'
$ws.Cells.Item(24, 10).Value = '[0.35757599 1.20549806 0.74883568]'
$ws.Cells.Item(24, 11).Value = 0.01736039522082607
$ws.Cells.Item(24, 12).Value = 0.7233400416045787
$ws.Cells.Item(25, 9).Value = '[-0.05784648 -0.0123828  -0.03122982]'
$ws.Cells.Item(25, 10).Value = '[0.53898832 1.63590089 0.95430879]'
$ws.Cells.Item(25, 11).Value = -0.03261981181292788
$ws.Cells.Item(25, 12).Value = 0.9814687983327414
$ws.Cells.Item(26, 9).Value = '[-0.01540889  0.04321879  0.00736465]'
$ws.Cells.Item(26, 10).Value = '[0.05780292 0.5934281  0.27955383]'
$ws.Cells.Item(26, 11).Value = 0.004404547690925986
$ws.Cells.Item(26, 12).Value = 0.361551570170175
$ws.Cells.Item(27, 9).Value = '[-0.04169891 -0.0050224  -0.02409862]'
$ws.Cells.Item(27, 10).Value = '[0.35129724 1.29457786 0.76425351]'
$ws.Cells.Item(27, 11).Value = -0.02744395360408734
$ws.Cells.Item(27, 12).Value = 0.820622702772395
$ws.Cells.Item(28, 9).Value = '[0.06156326 0.09779605 0.08039953]'
$ws.Cells.Item(28, 10).Value = '[1.82260039 2.64813738 2.3072696 ]'
$ws.Cells.Item(28, 11).Value = 0.06935041468073087
$ws.Cells.Item(28, 12).Value = 2.025612249610985
$ws.Cells.Item(29, 9).Value = '[-0.051475    0.00075582 -0.01915073]'
$ws.Cells.Item(29, 10).Value = '[0.22183313 1.22064741 0.63167042]'
$ws.Cells.Item(29, 11).Value = -0.01804290778053517
$ws.Cells.Item(29, 12).Value = 0.6244037625518721
$ws.Cells.Item(30, 9).Value = '[0.01860397 0.06487932 0.03581841]'
$ws.Cells.Item(30, 10).Value = '[0.67446195 1.5104682  1.18135252]'
$ws.Cells.Item(30, 11).Value = 0.02755623441823547
$ws.Cells.Item(30, 12).Value = 0.9262829919343325
$ws.Cells.Item(31, 9).Value = '[-0.02804289  0.01828391 -0.00659311]'
$ws.Cells.Item(31, 10).Value = '[0.04741339 0.54444848 0.28925011]'
$ws.Cells.Item(31, 11).Value = -0.008803736449432723
$ws.Cells.Item(31, 12).Value = 0.4288266963607174
$ws.Cells.Item(32, 9).Value = '[-0.03784661  0.00295076 -0.02003576]'
$ws.Cells.Item(32, 10).Value = '[0.195319   1.09051046 0.65884001]'
$ws.Cells.Item(32, 11).Value = -0.01939510909328935
$ws.Cells.Item(32, 12).Value = 0.6530236714976903
$ws.Cells.Item(33, 9).Value = '[-0.00600006  0.02690671  0.01022314]'
$ws.Cells.Item(33, 10).Value = '[0.26112234 1.07982612 0.55504504]'
$ws.Cells.Item(33, 11).Value = 0.002544975261400312
$ws.Cells.Item(33, 12).Value = 0.4604578166404271
$ws.Cells.Item(34, 9).Value = '[-0.04402136  0.01386575 -0.0184523 ]'
$ws.Cells.Item(34, 10).Value = '[0.07997376 1.10408141 0.39856426]'
$ws.Cells.Item(34, 11).Value = 0.0170061597529985
$ws.Cells.Item(34, 12).Value = 0.5024875591707831
$ws.Cells.Item(35, 9).Value = '[-0.03121757  0.00476695 -0.01581649]'
$ws.Cells.Item(35, 10).Value = '[0.23671329 0.90353526 0.47817032]'
$ws.Cells.Item(35, 11).Value = -0.01874019484844127
$ws.Cells.Item(35, 12).Value = 0.6490714767615496
$ws.Cells.Item(36, 9).Value = '[-0.037953    0.01088388 -0.01711748]'
$ws.Cells.Item(36, 10).Value = '[0.07176622 0.7509268  0.31604598]'
$ws.Cells.Item(36, 11).Value = -0.01306973580615289
$ws.Cells.Item(36, 12).Value = 0.4443302795451649
$ws.Cells.Item(37, 9).Value = '[0.02138275 0.07598825 0.05314268]'
$ws.Cells.Item(37, 10).Value = '[0.59942428 2.06766064 1.44650003]'
$ws.Cells.Item(37, 11).Value = 0.03894475938752369
$ws.Cells.Item(37, 12).Value = 1.130192666589308
$ws.Cells.Item(38, 9).Value = '[-0.01460162  0.03719121  0.01579296]'
$ws.Cells.Item(38, 10).Value = '[0.14210684 1.16773226 0.50139428]'
$ws.Cells.Item(38, 11).Value = 0.01559206625622982
$ws.Cells.Item(38, 12).Value = 0.5607397243416734
$ws.Cells.Item(39, 9).Value = '[-0.03054951  0.01481139 -0.0094538 ]'
$ws.Cells.Item(39, 10).Value = '[0.0834836  0.75742149 0.36478732]'
$ws.Cells.Item(39, 11).Value = 0.003552903535167937
$ws.Cells.Item(39, 12).Value = 0.4191440709826179
$ws.Cells.Item(40, 9).Value = '[0.0379805  0.0792435  0.05995563]'
$ws.Cells.Item(40, 10).Value = '[1.16795714 2.36544979 1.80692295]'
$ws.Cells.Item(40, 11).Value = 0.06609030530003528
$ws.Cells.Item(40, 12).Value = 1.910005731279309
$ws.Cells.Item(41, 9).Value = '[-0.01765365  0.0260773   0.00484041]'
$ws.Cells.Item(41, 10).Value = '[0.11295851 0.63308942 0.35205479]'
$ws.Cells.Item(41, 11).Value = 0.0010062841322327
$ws.Cells.Item(41, 12).Value = 0.4581731532815688
$ws.Cells.Item(42, 9).Value = '[-0.01734961  0.04364754  0.01976684]'
$ws.Cells.Item(42, 10).Value = '[0.2114115  1.28662873 0.57422533]'
$ws.Cells.Item(42, 11).Value = 0.008520923253518808
$ws.Cells.Item(42, 12).Value = 0.4610652262270156
$ws.Cells.Item(43, 9).Value = '[-0.032022    0.01005132 -0.01479994]'
$ws.Cells.Item(43, 10).Value = '[0.08349935 0.77016465 0.4091643 ]'
$ws.Cells.Item(43, 11).Value = -0.02107297226908885
$ws.Cells.Item(43, 12).Value = 0.6945123576849893
$ws.Cells.Item(44, 9).Value = '[0.01024687 0.06560455 0.03224365]'
$ws.Cells.Item(44, 10).Value = '[0.43812473 1.72671674 0.98009796]'
$ws.Cells.Item(44, 11).Value = 0.0160399383018276
$ws.Cells.Item(44, 12).Value = 0.6402461118553168
$ws.Cells.Item(45, 9).Value = '[-0.0326957   0.01344777 -0.01163711]'
$ws.Cells.Item(45, 10).Value = '[0.11159927 0.77090425 0.3724081 ]'
$ws.Cells.Item(45, 11).Value = -0.01635946191379515
$ws.Cells.Item(45, 12).Value = 0.582576571296237
$ws.Cells.Item(46, 9).Value = '[0.02010219 0.08079353 0.04627399]'
$ws.Cells.Item(46, 10).Value = '[0.5110349  2.15924686 1.20046219]'
$ws.Cells.Item(46, 11).Value = 0.0447784502840325
$ws.Cells.Item(46, 12).Value = 1.306973857846215
$ws.Cells.Item(47, 9).Value = '[0.0527318  0.11697764 0.07797485]'
$ws.Cells.Item(47, 10).Value = '[1.73526373 2.81952362 2.36673983]'
$ws.Cells.Item(47, 11).Value = 0.08226663966079767
$ws.Cells.Item(47, 12).Value = 2.341063992402707
$ws.Cells.Item(48, 9).Value = '[-0.01326739  0.04897825  0.01607394]'
$ws.Cells.Item(48, 10).Value = '[0.19767521 1.44202547 0.56225734]'
$ws.Cells.Item(48, 11).Value = 0.006416749046292192
$ws.Cells.Item(48, 12).Value = 0.3130425040857876
$ws.Cells.Item(49, 9).Value = '[0.07665479 0.13070618 0.10179546]'
$ws.Cells.Item(49, 10).Value = '[2.42526562 3.47101936 3.06300813]'
$ws.Cells.Item(49, 11).Value = 0.08763247384536613
$ws.Cells.Item(49, 12).Value = 2.495345542920472
$ws.Cells.Item(50, 9).Value = '[0.03210729 0.07640464 0.05738825]'
$ws.Cells.Item(50, 10).Value = '[0.87428856 2.22269542 1.58287777]'
$ws.Cells.Item(50, 11).Value = 0.04401790521592773
$ws.Cells.Item(50, 12).Value = 1.302491697519942
$ws.Cells.Item(51, 9).Value = '[0.01006408 0.06916056 0.03345107]'
$ws.Cells.Item(51, 10).Value = '[0.28332554 1.87687817 0.94682179]'
$ws.Cells.Item(51, 11).Value = -0.001645580655472509
$ws.Cells.Item(51, 12).Value = 0.2465649122155701
$ws.Cells.Item(52, 9).Value = '[0.00245986 0.05085734 0.02798822]'
$ws.Cells.Item(52, 10).Value = '[0.39977886 1.53219141 0.88130768]'
$ws.Cells.Item(52, 11).Value = 0.01425021219059619
$ws.Cells.Item(52, 12).Value = 0.5821229759745769
$ws.Cells.Item(53, 9).Value = '[-0.02999601  0.00665307 -0.0126243 ]'
$ws.Cells.Item(53, 10).Value = '[0.21330831 0.7921135  0.51089653]'
$ws.Cells.Item(53, 11).Value = -0.01932372502371776
$ws.Cells.Item(53, 12).Value = 0.7098045986669778
